$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells (H1:J1) + new data cells (H2:J2) -----------------
# Writing the values first registers the new shared strings in the order
# they appear in the target file (Usia, Gender, Asal Daerah).
$ws.Range("H1").Value = "Usia"
$ws.Range("I1").Value = "Gender"
$ws.Range("J1").Value = "Asal Daerah"

$ws.Range("H2").Value = 20
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1

# --- 2. Preserve the old "last column" look (border on the right) ---------
# G1/G2 used to be the last column and carried the right-hand border style;
# now J1/J2 are the new last column, so they should inherit that look
# before G1/G2 gets restyled to a "middle" column below.
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 3. G1/G2 become regular "middle" columns ------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 4. H2 ("Usia" value) gets its own small Arial font -------------------
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H2").Font.Size = 10
$ws.Range("H2").Font.Color = 2236704
$ws.Range("H2").Font.Name = "Arial"

# --- 5. Update dimension / selection ---------------------------------------
$ws.Range("A1:J2").Select() | Out-Null
